{"js": "// The edit removes the empty \"DACHUONG\"-styled paragraph that sits\n// directly after the \"5.2 Thi\u1ebft k\u1ebf gi\u1ea3i thu\u1eadt \u0111i\u1ec1u khi\u1ec3n khi\u1ec3n\" heading\n// paragraph (and immediately before the \"C\u00e1ch t\u01b0\u01a1ng t\u00e1c chung...\" body\n// paragraph), merging the two neighboring paragraphs together.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst headingText = \"5.2 Thi\u1ebft k\u1ebf gi\u1ea3i thu\u1eadt \u0111i\u1ec1u khi\u1ec3n khi\u1ec3n\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  const current = paragraphs.items[i];\n  const next = paragraphs.items[i + 1];\n  if (current.text.trim() === headingText && next.text.trim() === \"\") {\n    target = next;\n    break;\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# The edit removes the empty \"DACHUONG\"-styled paragraph that sits\n# directly after the \"5.2 Thi\u1ebft k\u1ebf gi\u1ea3i thu\u1eadt \u0111i\u1ec1u khi\u1ec3n khi\u1ec3n\" heading\n# paragraph (and immediately before the \"C\u00e1ch t\u01b0\u01a1ng t\u00e1c chung...\" body\n# paragraph), merging the two neighboring paragraphs together.\n\n$d = $word.ActiveDocument\n\n$headingText = \"5.2 Thi\u1ebft k\u1ebf gi\u1ea3i thu\u1eadt \u0111i\u1ec1u khi\u1ec3n khi\u1ec3n\"\n\n$headingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $headingText) {\n        $headingIndex = $i\n        break\n    }\n}\n\nif ($headingIndex -gt 0 -and $headingIndex -lt $d.Paragraphs.Count) {\n    $nextPara = $d.Paragraphs.Item($headingIndex + 1)\n    $nextText = $nextPara.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($nextText -eq \"\") {\n        $nextPara.Range.Delete()\n    }\n}\n"}
